$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Check")

$ws1.Range("A2").Value = 3016
$ws1.Range("E2").Value = 46200608016
$ws1.Range("X2").Value = "DN4127460130016"

$ws2.Range("A2").Value = 3016
$ws2.Range("C2").Value = "DN4127460130016"
